$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.245.38"
$ws.Range("E2").Value = "  -3.91%  "
$ws.Range("D3").Value = "2.609.36"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'518.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "'142.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "3.066.54"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("D14").Value = "58.207.89"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").Value = "'20.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "2.606.52"
$ws.Range("E17").Value = "  -8.77%  "
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").Value = "'335.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "'6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'64.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "0.0₃0791"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").Value = "'6.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'150.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "'4.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").Value = "'0.890"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").Value = "'0.849"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("D38").Value = "'36.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").Value = "'1.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.78%  "
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "'269.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "'19.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.05%  "
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "2.039.67"
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("D49").Value = "'0.0228"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'4.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("D51").Value = "'18.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
